$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.06401555221885813
$ws.Range("E2").Value = 0.06401555221885813

# Row 3
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.5578998287766188
$ws.Range("E3").Value = 0.5578998287766188

# Row 4
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.9638272997445375
$ws.Range("E4").Value = 0.9638272997445375

# Row 5
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = 0.5467079596503565
$ws.Range("E5").Value = 0.5467079596503565

# Row 6
$ws.Range("D6").Value = 0.06494869843063092
$ws.Range("E6").Value = 0.06494869843063092

# Row 7
$ws.Range("D7").Value = 0.261455180963484
$ws.Range("E7").Value = 0.738544819036516

# Row 8
$ws.Range("D8").Value = 0.001349037675607464
$ws.Range("E8").Value = 0.9986509623243925

# Row 9
$ws.Range("D9").Value = 0.2564839591102162
$ws.Range("E9").Value = 0.7435160408897838

# Row 10
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.544293350154333
$ws.Range("E10").Value = 0.455706649845667

# Row 11
$ws.Range("D11").Value = 0.1462554457326323
$ws.Range("E11").Value = 0.8537445542673676
$ws.Range("F11").Value = 1.690141081809998
$ws.Range("G11").Value = 0.3

# Row 12
$ws.Range("C12").Value = $true
$ws.Range("D12").Value = 0.0280263227034615
$ws.Range("E12").Value = 0.0280263227034615

# Row 13
$ws.Range("D13").Value = 0.07652751612780348
$ws.Range("E13").Value = 0.07652751612780348

# Row 14
$ws.Range("D14").Value = 0.9911137822305138
$ws.Range("E14").Value = 0.9911137822305138

# Row 15
$ws.Range("C15").Value = $false
$ws.Range("D15").Value = 0.7118829222330738
$ws.Range("E15").Value = 0.7118829222330738

# Row 16
$ws.Range("D16").Value = 0.02558194014012758
$ws.Range("E16").Value = 0.02558194014012758

# Row 17
$ws.Range("D17").Value = 0.2006394850985493
$ws.Range("E17").Value = 0.7993605149014507

# Row 18
$ws.Range("D18").Value = 0.00005710622041613605
$ws.Range("E18").Value = 0.9999428937795839

# Row 19
$ws.Range("D19").Value = 0.1533877439450247
$ws.Range("E19").Value = 0.8466122560549753

# Row 20
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = 0.6283050990847286
$ws.Range("E20").Value = 0.3716949009152714

# Row 21
$ws.Range("D21").Value = 0.06776519152352389
$ws.Range("E21").Value = 0.9322348084764761
$ws.Range("F21").Value = 2.250966310501099
$ws.Range("G21").Value = 0.4
